$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlPasteValues
$xlPasteValues = -4163

function Swap-Cell($ws, $r1, $c1, $r2, $c2, $scratchR, $scratchC) {
    $ws.Cells.Item($r1, $c1).Copy()
    $ws.Cells.Item($scratchR, $scratchC).PasteSpecial($xlPasteValues)

    $ws.Cells.Item($r2, $c2).Copy()
    $ws.Cells.Item($r1, $c1).PasteSpecial($xlPasteValues)

    $ws.Cells.Item($scratchR, $scratchC).Copy()
    $ws.Cells.Item($r2, $c2).PasteSpecial($xlPasteValues)
}

# Row 5 (EB-150KW/150KW 1410 x 1010 ...) and Row 6 (EB-180DC/120KB 1200 x 1600 ...)
# swap their Material / Quantity / Inventoryitem content (columns D, E, H).
Swap-Cell $ws 5 4 6 4 100 1
Swap-Cell $ws 5 5 6 5 100 2
Swap-Cell $ws 5 8 6 8 100 3

# Row 7 (EB-150KW/150KW 1010 X 1410 ...) and Row 8 (EB-180DC/120KB 1600 x 1200 ...)
# swap their Material / Quantity / Inventoryitem content (columns D, E, H).
Swap-Cell $ws 7 4 8 4 100 1
Swap-Cell $ws 7 5 8 5 100 2
Swap-Cell $ws 7 8 8 8 100 3

# Remove the scratch row used as temporary swap storage.
$ws.Rows.Item(100).Delete()
